$wb = $excel.ActiveWorkbook

# Sheet "runManager": update Execute flags from "no" to "yes" for several test cases
$wsRun = $wb.Worksheets.Item("runManager")
$wsRun.Range("B2").Value = "yes"
$wsRun.Range("B3").Value = "yes"
$wsRun.Range("B6").Value = "yes"
$wsRun.Range("B7").Value = "yes"

# Sheet "iterationdata": becomes the active sheet, selection moves to B4
$wsIter = $wb.Worksheets.Item("iterationdata")
$wsIter.Activate()
$wsIter.Range("B4").Select()
